$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "URL"
$ws.Range("D35").Select()
